$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E1 from a numeric date value to the text label "03_03_2024"
# (matches the other header cells B1/C1/D1 which are plain text date labels)
$ws.Range("E1").Value = "03_03_2024"

# Update the active selection to E2 (as recorded in the saved file)
$ws.Range("E2").Select()
